# daily auto push: 2026-01-13 13:50 UTC
# A new observation (2026/01/13, 火, 20, 29) is inserted into the log at
# row 617, pushing all later rows down by one (617->618, ..., 658->659).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 617; every row from 617 on shifts down by one
# and the sheet grows from A1:D658 to A1:D659 automatically.
$ws.Rows.Item(617).Insert()

# Fill the newly inserted row. The date column stores plain text like
# "2026/01/13" (not an Excel date), so prefix with an apostrophe to force
# text entry, then clear the resulting formatting so no stray style
# (quote-prefix / number format) is left behind on the cell.
$ws.Range("A617").Value = "'2026/01/13"
$ws.Range("A617").ClearFormats()
$ws.Range("B617").Value = "火"
$ws.Range("C617").Value = 20
$ws.Range("D617").Value = 29
